$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.016.49'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '3.421.03'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'409.89"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').Value = "'128.86"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.16%  '
$ws.Range('D7').Value = "'0.636"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +7.69%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +6.00%  '
$ws.Range('D10').Value = "'0.140"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +11.68%  '
$ws.Range('D11').Value = "'42.77"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'0.141"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'9.07"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +7.63%  '
$ws.Range('D14').Value = '3.956.99'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = "'21.27"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +7.37%  '
$ws.Range('D16').Value = "'0.0000208"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +48.79%  '
$ws.Range('D17').Value = '3.412.96'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = "'12.25"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +4.49%  '
$ws.Range('E19').Value = '  +7.16%  '
$ws.Range('D20').Value = '61.960.20'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = "'447.33"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +43.52%  '
$ws.Range('D22').Value = "'92.07"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +10.42%  '
$ws.Range('D23').Value = "'3.17"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('E25').Value = '  +3.06%  '
$ws.Range('D26').Value = "'33.43"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +12.95%  '
$ws.Range('D27').Value = "'8.88"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +10.38%  '
$ws.Range('D28').Value = "'4.78"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').Value = "'2.75"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D30').Value = "'7.57"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -4.87%  '
$ws.Range('D31').Value = "'12.01"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +5.85%  '
$ws.Range('D32').Value = "'0.169"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = "'42.91"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = "'0.114"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = "'0.0498"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.34%  '
$ws.Range('D37').Value = "'53.45"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +4.09%  '
$ws.Range('D38').Value = "'0.997"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').Value = "'3.38"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('E40').Value = '  +7.56%  '
$ws.Range('E41').Value = '  -0.90%  '
$ws.Range('D42').Value = "'0.313"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.98%  '
$ws.Range('D43').Value = "'141.41"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.19%  '
$ws.Range('D44').Value = "'4.23"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +6.66%  '
$ws.Range('D45').Value = "'1.98"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('E46').Value = '  +8.64%  '
$ws.Range('D47').Value = "'16.57"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = "'22.39"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +5.11%  '
$ws.Range('D49').Value = "'2.12"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +8.46%  '
$ws.Range('D50').Value = '3.762.64'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').Value = '2.116.36'
$ws.Range('E51').Value = '  +0.45%  '
